# Weekly update: a new week of data (rows 22-23) is inserted at the top of the
# "shifting window" (rows 22-69). Every row from 24..69 takes on the values
# that used to live two rows above it (r-2), the two oldest rows (old 68, 69)
# overflow into two brand-new rows (70, 71), and rows 22-23 get the new
# week's own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 22
$lastRowBefore = 69
$numCols = 18   # columns A..R

# 1) Snapshot the current contents (values only) of rows 22..69 before we
#    start overwriting anything, so the shift doesn't clobber data it still
#    needs to read from.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRowBefore; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2()
    }
    $snapshot[$r] = $rowVals
}

# Reusable date format, taken from the existing D column (style s="2").
$dateFmt = $ws.Cells.Item($firstRow, 4).NumberFormat()

# 2) Shift rows 24..69 down from what used to be two rows above (old 22..67),
#    then spill the last two surviving rows (old 68, 69) into new rows 70, 71.
for ($r = $lastRowBefore; $r -ge ($firstRow + 2); $r--) {
    $src = $snapshot[$r - 2]
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($r, $c).Value = $src[$c - 1]
    }
}

$newRow1 = $lastRowBefore + 1   # 70
$newRow2 = $lastRowBefore + 2   # 71

$src = $snapshot[$lastRowBefore - 1]   # old row 68
for ($c = 1; $c -le $numCols; $c++) {
    $ws.Cells.Item($newRow1, $c).Value = $src[$c - 1]
}
$ws.Cells.Item($newRow1, 4).NumberFormat = $dateFmt

$src = $snapshot[$lastRowBefore]       # old row 69
for ($c = 1; $c -le $numCols; $c++) {
    $ws.Cells.Item($newRow2, $c).Value = $src[$c - 1]
}
$ws.Cells.Item($newRow2, 4).NumberFormat = $dateFmt

# 3) Rows 22 and 23 hold the brand-new week's prices: same market/category/
#    variety/quality/unit metadata as before, but a later date, a new
#    volume figure and a different origin province.
$ws.Cells.Item(22, 4).Value = 44519            # Fecha
$ws.Cells.Item(22, 10).Value = 270             # Volumen
$ws.Cells.Item(22, 15).Value = "Provincia de Linares"   # Origen

$ws.Cells.Item(23, 4).Value = 44519            # Fecha
$ws.Cells.Item(23, 10).Value = 380             # Volumen
$ws.Cells.Item(23, 15).Value = "Provincia de Linares"   # Origen
